# [Feat : KSW] Add Monster Table, DataManager SingleTon
#
# Adds a second row to the Enum sheet describing a new "MonsterGrade"
# category with a "Common"/"Boss" pair, widens column A to fit the new
# header text, moves the active selection to C3, and sets the page setup
# (paper size / orientation) that Excel stamps once a sheet has been
# printed/configured.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row of data: A2="MonsterGrade", B2="Common" (reuses existing shared
# string), C2="Boss"
$ws.Range("A2").Value = "MonsterGrade"
$ws.Range("B2").Value = "Common"
$ws.Range("C2").Value = "Boss"

# Column A needs to be wide enough to show "MonsterGrade" in full.
$ws.Columns.Item(1).ColumnWidth = 12

# Selection moves to C3 after the edits.
$ws.Range("C3").Select() | Out-Null

# Page setup: Letter/A4-adjacent "paperSize=9" (A4) in portrait orientation.
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1
